$d = $word.ActiveDocument

# Remove the "010782915" paragraph and the empty centered paragraph that
# follows it, leaving the "Zi Yan Zhang" paragraph directly followed by the
# write-up body paragraph.
$start = $d.Paragraphs(2).Range.Start
$end = $d.Paragraphs(3).Range.End
$r = $d.Range($start, $end)
$r.Delete()
